# Change the table style (table design) of the three data tables on
# slides 14, 15 and 16 from the default "Table_0" medium-style
# ({C3C84C1C-FEA8-4A95-AE0F-EB79AEA6F25F}) to the built-in
# "Medium Style 2 - Accent 1" style ({1F88DE88-8345-425C-85AB-6B369D91272E}).
#
# Table styles can't be assigned through the `Table.Style` property —
# PowerPoint requires `Table.ApplyStyle("{GUID}")` — so we locate the
# single table shape on each of those slides and re-apply the style via
# that method.

$p = $ppt.ActivePresentation

$newStyleId = "{1F88DE88-8345-425C-85AB-6B369D91272E}"
$slideNumbers = 14, 15, 16

foreach ($slideNumber in $slideNumbers) {
    $slide = $p.Slides.Item($slideNumber)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
